# Update the "2024" worksheet: a new September entry was logged
# ("bal axisbank w axis" at 2024-09-04 06:53:15). This is inserted as a
# brand-new row 29, pushing all the existing rows 29-51 down by one
# (29->30, 30->31, ..., 51->52), which is exactly the shift pattern shown
# in the diff (each row's September Details/Date move to the row below,
# the August hdfc chain shifts the same way once the September chain runs
# out, and the trailing "Broadband" label ends up on the new last row, 52).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a fresh blank row at row 29; everything below shifts down one row.
$ws.Rows.Item(29).Insert()

# Populate the new row's September Details / September Date cells.
$ws.Range("R29").Value = "bal axisbank w axis"
$ws.Range("S29").Value = "2024-09-04 06:53:15"
